# Cronograma BD.xlsx - "Conexion a BD y correccion de cronograma"
# Corrects the schedule (Gantt) start date and several milestone rows,
# and updates the sheet view (zoom/selection) to match the author's session.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Gantt")

# --- Project start date moved 2 days later (Inicio_del_proyecto, F3) ---
$ws.Range("F3").Value = 43699

# --- Milestone rows (table "Hitos", B10:G14) -----------------------------
# Row 10: PERSPECTVA DEL ÁREA DE BASE DE DATOS
$ws.Range("E10").Value = 0.5
$ws.Range("G10").Value = 5

# Row 11: MODELOS DE DATOS
$ws.Range("F11").Value = 43705
$ws.Range("G11").Value = 5

# Row 12: MODELO ENTIDAD - RELACIÓN
$ws.Range("E12").Value = 0.2
$ws.Range("F12").Value = 43682
$ws.Range("G12").Value = 5

# Row 13: DISEÑO DE BASES DE DATOS RELACIONALES
$ws.Range("E13").Value = 0.2
$ws.Range("F13").Value = 43713
$ws.Range("G13").Value = 5

# Row 14: RECUPERACIÓN DE FALLAS
$ws.Range("F14").Value = 43717
$ws.Range("G14").Value = 5

# --- Sheet view: zoom out and move the selection ---------------------------
$ws.Activate()
$ws.Range("E17").Select()
$excel.ActiveWindow.Zoom = 55
